{"js": "// Replace every arithmetic-expression cell in the worksheet's single table\n// with the corresponding new expression from the commit's diff. The table\n// is 20 rows x 5 columns (100 cells); values below are listed in row-major\n// order (same order the cells appear in the document / diff).\nconst newValues = [\n  [\"46+0=\", \"50+37=\", \"13+19=\", \"96-20=\", \"65+10=\"],\n  [\"27+60=\", \"9+32=\", \"85-20=\", \"72+23=\", \"25+7=\"],\n  [\"75-40=\", \"56-52=\", \"44+52=\", \"18+5=\", \"34+41=\"],\n  [\"41+47=\", \"35+33=\", \"7+65=\", \"4+45=\", \"76-72=\"],\n  [\"26+5=\", \"79-68=\", \"24+62=\", \"52+3=\", \"98-89=\"],\n  [\"80-60=\", \"55-5=\", \"96-71=\", \"54-31=\", \"64-13=\"],\n  [\"17+29=\", \"4+72=\", \"94-18=\", \"67-52=\", \"40-12=\"],\n  [\"99-47=\", \"36-27=\", \"40+4=\", \"64-31=\", \"67+5=\"],\n  [\"82-52=\", \"83+6=\", \"74+4=\", \"47+19=\", \"80-64=\"],\n  [\"68+28=\", \"8+30=\", \"79-64=\", \"55-10=\", \"51-9=\"],\n  [\"26-8=\", \"88-56=\", \"51+22=\", \"86-29=\", \"71+9=\"],\n  [\"67+25=\", \"89-28=\", \"46+9=\", \"43+17=\", \"58-22=\"],\n  [\"92-10=\", \"7+27=\", \"84-9=\", \"44-34=\", \"7-6=\"],\n  [\"23-3=\", \"90-54=\", \"30+49=\", \"81-55=\", \"5+59=\"],\n  [\"6+12=\", \"3-3=\", \"27+10=\", \"29+9=\", \"26+9=\"],\n  [\"87-80=\", \"71-49=\", \"68+19=\", \"38-29=\", \"60+15=\"],\n  [\"17+28=\", \"14+68=\", \"84-68=\", \"40+24=\", \"84-8=\"],\n  [\"89-80=\", \"29-14=\", \"35-7=\", \"14+70=\", \"28+67=\"],\n  [\"93-35=\", \"36+11=\", \"25+52=\", \"24+45=\", \"1+92=\"],\n  [\"1+90=\", \"20-5=\", \"78+12=\", \"8+34=\", \"99-33=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nif (table.rowCount === newValues.length) {\n  // Word.Table.values is a 2-D array projection of the cell text; assigning\n  // it rewrites each cell's Range.Text in place while leaving paragraph /\n  // run formatting (fonts, size, alignment) untouched.\n  table.values = newValues;\n  await context.sync();\n}\n", "ps1": "# Replace every arithmetic-expression cell in the worksheet's single table\n# with the corresponding new expression from the commit's diff. The table\n# is 20 rows x 5 columns (100 cells); values below are listed in row-major\n# order (same order the cells appear in the document / diff).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$values = @(\n  @(\"46+0=\",\"50+37=\",\"13+19=\",\"96-20=\",\"65+10=\"),\n  @(\"27+60=\",\"9+32=\",\"85-20=\",\"72+23=\",\"25+7=\"),\n  @(\"75-40=\",\"56-52=\",\"44+52=\",\"18+5=\",\"34+41=\"),\n  @(\"41+47=\",\"35+33=\",\"7+65=\",\"4+45=\",\"76-72=\"),\n  @(\"26+5=\",\"79-68=\",\"24+62=\",\"52+3=\",\"98-89=\"),\n  @(\"80-60=\",\"55-5=\",\"96-71=\",\"54-31=\",\"64-13=\"),\n  @(\"17+29=\",\"4+72=\",\"94-18=\",\"67-52=\",\"40-12=\"),\n  @(\"99-47=\",\"36-27=\",\"40+4=\",\"64-31=\",\"67+5=\"),\n  @(\"82-52=\",\"83+6=\",\"74+4=\",\"47+19=\",\"80-64=\"),\n  @(\"68+28=\",\"8+30=\",\"79-64=\",\"55-10=\",\"51-9=\"),\n  @(\"26-8=\",\"88-56=\",\"51+22=\",\"86-29=\",\"71+9=\"),\n  @(\"67+25=\",\"89-28=\",\"46+9=\",\"43+17=\",\"58-22=\"),\n  @(\"92-10=\",\"7+27=\",\"84-9=\",\"44-34=\",\"7-6=\"),\n  @(\"23-3=\",\"90-54=\",\"30+49=\",\"81-55=\",\"5+59=\"),\n  @(\"6+12=\",\"3-3=\",\"27+10=\",\"29+9=\",\"26+9=\"),\n  @(\"87-80=\",\"71-49=\",\"68+19=\",\"38-29=\",\"60+15=\"),\n  @(\"17+28=\",\"14+68=\",\"84-68=\",\"40+24=\",\"84-8=\"),\n  @(\"89-80=\",\"29-14=\",\"35-7=\",\"14+70=\",\"28+67=\"),\n  @(\"93-35=\",\"36+11=\",\"25+52=\",\"24+45=\",\"1+92=\"),\n  @(\"1+90=\",\"20-5=\",\"78+12=\",\"8+34=\",\"99-33=\")\n)\n\nif ($tbl.Rows.Count -eq $values.Count) {\n  for ($r = 1; $r -le $values.Count; $r++) {\n    $row = $values[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n      # Setting Range.Text replaces only the run text, preserving the\n      # paragraph/run formatting (fonts, size, alignment) already on the cell.\n      $tbl.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n  }\n}"}
